# Append new scrape batch (2025-09-20 06:23:53 JST) to the "ランサーズ" sheet.
# The oldest rows fall off the bottom (sheet is capped at 6 data rows + header),
# newest rows are written starting at row 2, everything shifts down no longer
# than the cap, and the trailing rows (previously 8-12) are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (2025-09-20 06:23:53 batch), in final on-sheet order (row 2..7).
$timestamp = "2025-09-20 06:23:53"

$rows = @(
    @{ B = "【急募】スマホアプリ自動化デモ開発(LLM連携)"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397035"; G = 228; H = "★スマホアプリ ◆開発,自動化 ◇アプリ" },
    @{ B = "システム開発において活躍できる案件紹介"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397117"; G = 125; H = "◆開発,システム開発" },
    @{ B = "システム開発の案件紹介"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397121"; G = 125; H = "◆開発,システム開発" },
    @{ B = "システム開発の複数案件紹介"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397127"; G = 125; H = "◆開発,システム開発" },
    @{ B = "【フォートナイト】クリエイティブ作品を世界に公開したい!"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397192"; G = 18; H = $null },
    @{ B = "初回 Web広告のタグ設置・動作確認"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397007"; G = 18; H = $null }
)

# Drop every existing hyperlink up front -- the sheet is about to be rewritten
# wholesale and stale relationships (rows 8-12) must not survive.
$ws.Range("A1").Hyperlinks.Delete()

$lastRow = 1 + $rows.Count

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Value = $data.F
    $ws.Hyperlinks.Add($fCell, $data.F)

    $ws.Cells.Item($r, 7).Value = $data.G

    if ($null -eq $data.H) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $data.H
    }
}

# Remove now-stale rows that used to hold the older (now dropped) entries.
$ws.Rows("8:12").Delete()

# Column widths: B and D both narrow from their old widths down to 30 chars.
# (ColumnWidth is character-width; the saved OOXML <col width> ends up 5/6
# wider, so dial the COM value back by 5/6 to land on an exact 30.)
$ws.Columns.Item(2).ColumnWidth = 30 - 5/6
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
